# Add a new data row (row 8) to the sheet, mirroring the existing rows'
# layout: machine name (text), temperature (number), pression (number),
# and an empty timestamp cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "mskdgjmsldfgmdsfgjdsmfg"
$ws.Range("B8").Value = 2342354
$ws.Range("C8").Value = 324252346.5
